# Update countries & provincias Spain
# Applies the 24-Jul-2020 00:55 data refresh to the "Pais" sheet:
#  - refreshed case counts for a number of countries (numbers only)
#  - three country pairs that swap ranking order because one of the pair
#    received updated (higher) numbers that now outrank its neighbour:
#      Turquia / Colombia           (rows 18/19)
#      Jamaica / Togo                (rows 151/152)
#      Groenlandia / Islas Malvinas  (rows 210/211, tied counts)
#  - the "Datos actualizados..." timestamp banner in A1

function Set-Row8($Row, $A, $B, $C, $D, $E, $F, $G, $H) {
    $arr = New-Object 'object[,]' 1,8
    $arr[0,0] = $A
    $arr[0,1] = $B
    $arr[0,2] = $C
    $arr[0,3] = $D
    $arr[0,4] = $E
    $arr[0,5] = $F
    $arr[0,6] = $G
    $arr[0,7] = $H
    $ws.Range("A" + $Row + ":H" + $Row).Value = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header banner timestamp
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 00:55"

# Estados Unidos (row 4) - refreshed numbers
Set-Row8 4 "Estados Unidos" 4165113 64238 1977400 2040473 0 1057 147240

# Peru (row 9) - refreshed numbers
Set-Row8 9 "Peru" 371096 4546 255945 97497 0 199 17654

# Colombia overtakes Turquia: row 18 becomes Colombia (new numbers),
# row 19 becomes Turquia (its prior, unchanged numbers)
Set-Row8 18 "Colombia" 226373 7945 107951 110734 0 315 7688
Set-Row8 19 "Turquia"  223315 913  206365 11387  0 18  5563

# Argentina (row 23) - refreshed numbers
Set-Row8 23 "Argentina" 148027 6127 62815 82510 0 114 2702

# Canada (row 24) - refreshed numbers
Set-Row8 24 "Canada" 112672 432 98519 5279 0 4 8874

# Nigeria (row 50) - refreshed numbers
Set-Row8 50 "Nigeria" 38948 604 16061 22054 0 20 833

# Japon (row 59) - refreshed numbers
Set-Row8 59 "Japon" 27029 726 21035 5004 0 1 990

# Chequia (row 71) - refreshed numbers
Set-Row8 71 "Chequia" 14800 230 9328 5107 0 1 365

# Noruega (row 85) - refreshed numbers
Set-Row8 85 "Noruega" 9085 26 8674 156 0 0 255

# Guinea (row 93) - refreshed numbers
Set-Row8 93 "Guinea" 6806 59 5999 765 0 1 42

# Republica de Africa Central (row 99) - refreshed numbers
Set-Row8 99 "Republica de Africa Central" 4590 16 1452 3080 0 1 58

# Cabo Verde (row 121) - refreshed numbers
Set-Row8 121 "Cabo Verde" 2190 36 1150 1019 0 0 21

# Niger (row 141) - refreshed numbers
Set-Row8 141 "Niger" 1124 2 1022 33 0 0 69

# Togo overtakes Jamaica: row 151 becomes Togo (new numbers),
# row 152 becomes Jamaica (its prior, unchanged numbers)
Set-Row8 151 "Togo"    828 22 584 228 0 0 16
Set-Row8 152 "Jamaica" 816 6  710 96  0 0 10

# Islas Caimanes (row 174) - refreshed numbers
Set-Row8 174 "Islas Caimanes" 203 0 202 0 0 0 1

# Groenlandia / Islas Malvinas are tied (13/0/13/0/0/0/0); their relative
# order flips in the refreshed ranking, so swap the country names only
Set-Row8 210 "Islas Malvinas" 13 0 13 0 0 0 0
Set-Row8 211 "Groenlandia"    13 0 13 0 0 0 0
